$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B:E and G for rows 2-6 (F is unchanged)
$data = @{
    2 = @(0.0006075818656279264, 0.002658071450198252, 18.71679738969934, 2797.565817734744, 2816.285880777759)
    3 = @(0.01253208636536152, 0.3048912486333797, 3.223369029078222, 13.86384647080068, 17.40463883487765)
    4 = @(1.445647641019636, 1.626987699542094, 0.7210945179870265, 13.86384647080068, 17.65757632934944)
    5 = @(1.445647641019636, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 6.82939032824165)
    6 = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 8.656069925401464)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]  # B
    $ws.Cells.Item($row, 3).Value = $vals[1]  # C
    $ws.Cells.Item($row, 4).Value = $vals[2]  # D
    $ws.Cells.Item($row, 5).Value = $vals[3]  # E
    $ws.Cells.Item($row, 7).Value = $vals[4]  # G
}
